# Timesheet update: log 2 hours worked on Monday (row 18) and
# Tuesday (row 19) of the week, which rolls up into the weekly
# total (C24) and the grand total (C46) via the existing SUM formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 2

# Move the active selection to the next unfilled day, as the
# author did after making the entry.
[void]$ws.Range("C20").Select()
